# Apply the INS final-commit edits to Sheet1 of the BreastCancer INS workbook.
# The workbook already has the 5 SQL-query cells (B2, C2, B3, B4, B5) populated;
# this script rewrites their text to the corrected/updated query text, re-applies
# a uniform font/wrap format across them (mirrors the consolidated style Excel
# produces when the same formatting is (re)applied to the whole block in one go),
# and updates the current selection to reflect where the user ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits -------------------------------------------------------

# B2: "Programs" tab query - relabel "Focus Area" -> "Special Topic", tidy up
# the FROM/WHERE split and stray blank line.
$ws.Range("B2").Value = 'SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Special Topic",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.program_acronym     
        ELSE prg.data_link
    END AS "Data Location Details" 
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE ''%Breast Cancer%''
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;'

# C2: summary counts query - terminate the statement with a semicolon.
$ws.Range("C2").Value = 'SELECT DISTINCT
    COUNT(DISTINCT prg.program_id) AS "Programs",
    COUNT(DISTINCT prj.project_id) AS "Projects",
    COUNT(DISTINCT gnt.grant_id) AS "Grants",
    COUNT(DISTINCT pub.pmid) AS "Publications"
FROM 
    df_program prg
LEFT JOIN 
    df_project prj ON prg.program_id = prj."program.program_id"
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type LIKE ''%Breast Cancer%'';'

# B3: "Projects" tab query - org_name -> project_org_name.
$ws.Range("B3").Value = 'SELECT DISTINCT
    prj.project_id AS "Project ID", 
    prj.project_title AS "Project Title",
    prj.project_org_name AS "Organization",
    prj.project_start_date AS "Project Start Date",
    prj.project_end_date AS "Project End Date"
FROM 
    df_project prj
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
     prg.cancer_type LIKE ''%Breast Cancer%''
ORDER BY 
    lower(prj.project_id) ASC
LIMIT 100;'

# B4: "Grants" tab query - project_end_date -> grant_end_date, plus a stray
# double space before LIKE.
$ws.Range("B4").Value = 'SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.grant_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type  LIKE ''%Breast Cancer%''
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;'

# B5: "Publications" tab query - title -> publication_title, add a CASE
# branch for relative_citation_ratio = 1.0, plus a stray double space
# before LIKE.
$ws.Range("B5").Value = 'SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.publication_title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN ''0''
    WHEN pub.relative_citation_ratio = 7.0 THEN ''7''
    WHEN pub.relative_citation_ratio = 2.0 THEN ''2''
  WHEN pub.relative_citation_ratio = 1.0 THEN ''1''
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
     prg.cancer_type  LIKE ''%Breast Cancer%''
ORDER BY 
    lower(pub.pmid) ASC
LIMIT 100;'

# --- Formatting ------------------------------------------------------------
# Re-apply the same font/wrap formatting across all the query cells in a
# single pass so they share one consolidated style.
$fmtRange = $ws.Range("B2:C2,B3:B3,B4:B4,B5:B5")
$fmtRange.Font.Size = 12
$fmtRange.WrapText = $true

# --- View state --------------------------------------------------------
# Reflect the user's final selection/scroll position.
$ws.Range("D2").Select()
